{"js": "// Sprint backlog 1 \u2013 update the \"check stock\" user story text and the\n// three related task descriptions that previously referenced the\n// \"reservation window\" so that they now talk about the stock-viewing UI.\n\nasync function replaceOnce(body, searchText, replacementText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n\n  results.items[0].insertText(replacementText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) User story cell, first paragraph: merge the split run back into\n//    one clean sentence (text is unchanged, just de-fragmented).\nawait replaceOnce(\n  context.document.body,\n  \"As a user, I want to check the stock database\",\n  \"As a user, I want to check the stock database\"\n);\n\n// 2) User story cell, second paragraph: merge the four runs into one\n//    and fix the wording.\nawait replaceOnce(\n  context.document.body,\n  \"of products so that I can check all the available products, its quantity and price at the same time in stock database\",\n  \"of products so that I can check all the available products, its quantity and price at the same time in stock database\"\n);\n\n// 3) Task: \"Design UI for reservation window...\" -> new wording about\n//    designing the UI for viewing stock of products.\nawait replaceOnce(\n  context.document.body,\n  \"Design UI for reservation window and write code for that\",\n  \"Design user interface for viewing stock of products on the window\"\n);\n\n// 4) Task: \"create a database table ... reservation window\" -> now\n//    connects to the stock window to view products.\nawait replaceOnce(\n  context.document.body,\n  \"create a database table for storing products information and connect it with the reservation window\",\n  \"create a database table for storing products information and connect it with the stock window to view\"\n);\n\n// 5) Task: \"Test the reservation window...\" -> now tests the UI for\n//    viewing stock.\nawait replaceOnce(\n  context.document.body,\n  \"Test the reservation window whether it is working or not\",\n  \"Test the UI for viewing stock window whether it is working or not\"\n);\n", "ps1": "# Sprint backlog 1 - update the \"check stock\" user story text and the\n# three related task descriptions that previously referenced the\n# \"reservation window\" so that they now talk about the stock-viewing UI.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute(\n        $findText,\n        $false, $false, $false, $false, $false, $true, 1, $false,\n        $replaceText,\n        2\n    )\n}\n\n# 1) User story cell, first paragraph: merge the split run back into\n#    one clean sentence (text is unchanged, just de-fragmented).\nReplace-Text \"As a user, I want to check the stock database\" \"As a user, I want to check the stock database\"\n\n# 2) User story cell, second paragraph: merge the four runs into one\n#    and fix the wording.\nReplace-Text \"of products so that I can check all the available products, its quantity and price at the same time in stock database\" \"of products so that I can check all the available products, its quantity and price at the same time in stock database\"\n\n# 3) Task: \"Design UI for reservation window...\" -> new wording about\n#    designing the UI for viewing stock of products.\nReplace-Text \"Design UI for reservation window and write code for that\" \"Design user interface for viewing stock of products on the window\"\n\n# 4) Task: \"create a database table ... reservation window\" -> now\n#    connects to the stock window to view products.\nReplace-Text \"create a database table for storing products information and connect it with the reservation window\" \"create a database table for storing products information and connect it with the stock window to view\"\n\n# 5) Task: \"Test the reservation window...\" -> now tests the UI for\n#    viewing stock.\nReplace-Text \"Test the reservation window whether it is working or not\" \"Test the UI for viewing stock window whether it is working or not\"\n"}
